$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force B2:E51 to text format before writing, so that numeric-looking
# strings (e.g. "0.999", "605.21") are written as literal text rather than
# being auto-converted to numbers by Excel, matching the source inlineStr cells.
$ws.Range("B2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '66.117.27'
$ws.Range("E2").Value = '  +0.22%  '

$ws.Range("D3").Value = '3.171.69'
$ws.Range("E3").Value = '  -1.45%  '

$ws.Range("E4").Value = '  +0.05%  '

$ws.Range("D5").Value = '605.21'
$ws.Range("E5").Value = '  +0.10%  '

$ws.Range("D6").Value = '154.26'
$ws.Range("E6").Value = '  +0.19%  '

$ws.Range("D7").Value = '0.999'
$ws.Range("E7").Value = '  +0.00%  '

$ws.Range("D8").Value = '3.168.38'
$ws.Range("E8").Value = '  -1.51%  '

$ws.Range("E9").Value = '  +2.30%  '

$ws.Range("E10").Value = '  -1.20%  '

$ws.Range("D11").Value = '5.66'
$ws.Range("E11").Value = '  -8.06%  '

$ws.Range("D12").Value = '0.519'
$ws.Range("E12").Value = '  +1.70%  '

$ws.Range("D13").Value = '0.0000268'
$ws.Range("E13").Value = '  -1.23%  '

$ws.Range("D14").Value = '38.35'
$ws.Range("E14").Value = '  -2.05%  '

$ws.Range("D15").Value = '3.690.85'
$ws.Range("E15").Value = '  -1.49%  '

$ws.Range("D16").Value = '66.147.09'

$ws.Range("D17").Value = '7.43'
$ws.Range("E17").Value = '  -1.17%  '

$ws.Range("D18").Value = '3.169.28'
$ws.Range("E18").Value = '  -1.46%  '

$ws.Range("E19").Value = '  +0.90%  '

$ws.Range("D20").Value = '509.71'
$ws.Range("E20").Value = '  -0.16%  '

$ws.Range("D21").Value = '15.39'
$ws.Range("E21").Value = '  -0.69%  '

$ws.Range("D22").Value = '0.730'
$ws.Range("E22").Value = '  -1.14%  '

$ws.Range("D23").Value = '8.02'
$ws.Range("E23").Value = '  -0.91%  '

$ws.Range("D24").Value = '14.89'
$ws.Range("E24").Value = '  -3.23%  '

$ws.Range("D25").Value = '84.59'
$ws.Range("E25").Value = '  -0.58%  '

$ws.Range("E26").Value = '  +0.14%  '

$ws.Range("D27").Value = '3.00'
$ws.Range("E27").Value = '  -0.58%  '

$ws.Range("D28").Value = '9.15'
$ws.Range("E28").Value = '  -0.50%  '

$ws.Range("D29").Value = '2.38'
$ws.Range("E29").Value = '  +4.74%  '

$ws.Range("E30").Value = '  +5.19%  '

$ws.Range("D31").Value = '7.18'
$ws.Range("E31").Value = '  +5.33%  '

$ws.Range("D32").Value = '27.95'
$ws.Range("E32").Value = '  -0.93%  '

$ws.Range("E33").Value = '  +0.11%  '

$ws.Range("D34").Value = '1.20'
$ws.Range("E34").Value = '  -1.45%  '

$ws.Range("D35").Value = '6.52'
$ws.Range("E35").Value = '  -1.22%  '

$ws.Range("D36").Value = '502.48'
$ws.Range("E36").Value = '  +4.50%  '

$ws.Range("D37").Value = '54.98'
$ws.Range("E37").Value = '  -0.20%  '

$ws.Range("D38").Value = '0.0883'
$ws.Range("E38").Value = '  -2.30%  '

$ws.Range("D39").Value = '0.0419'
$ws.Range("E39").Value = '  -0.30%  '

$ws.Range("E40").Value = '  +6.52%  '

$ws.Range("D41").Value = '8.78'
$ws.Range("E41").Value = '  -1.60%  '

$ws.Range("D42").Value = '0.0₃0685'
$ws.Range("E42").Value = '  +6.52%  '

$ws.Range("D43").Value = '0.299'
$ws.Range("E43").Value = '  -0.16%  '

$ws.Range("D44").Value = '2.81'
$ws.Range("E44").Value = '  -4.56%  '

$ws.Range("D45").Value = '2.45'
$ws.Range("E45").Value = '  +0.07%  '

$ws.Range("D46").Value = '2.822.54'
$ws.Range("E46").Value = '  -4.27%  '

$ws.Range("D47").Value = '28.06'
$ws.Range("E47").Value = '  -1.97%  '

$ws.Range("B48").Value = 'ThetaToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D48").Value = '2.38'
$ws.Range("E48").Value = '  +3.11%  '

$ws.Range("B49").Value = 'USDe'
$ws.Range("C49").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D49").Value = '0.999'
$ws.Range("E49").Value = '  -0.11%  '

$ws.Range("D50").Value = '0.117'
$ws.Range("E50").Value = '  +0.49%  '

$ws.Range("D51").Value = '35.33'
$ws.Range("E51").Value = '  +6.82%  '

# Clear the temporary formatting again so the cell style indices are left
# as they were (no stray number-format style lingers on the written cells).
$ws.Range("B2:E51").ClearFormats()
